# Updated cryptos list (price + volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.257.67"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'1.850.19"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'241.26"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'0.6735"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07444"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "'0.2968"
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "'22.95"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "'1.819.14"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "'5.023"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'0.6797"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "'6.177"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'29.221.84"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'0.000008307"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'229.24"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").Value = "'12.58"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'7.232"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'160.83"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'8.703"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1416"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").Value = "'18.06"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'1.509"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "'4.194"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "'4.082"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").Value = "'1.187"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'0.05321"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").Value = "'0.7580"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "'1.141"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "'2.688"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "'1.335.42"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("D38").Value = "'0.01806"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "'2.738"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "'0.9254"
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").Value = "'5.964"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'103.64"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'1.978.50"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'0.07828"
$ws.Range("E45").Value = "  +7.17%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "'63.84"
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").Value = "'9.252"
$ws.Range("D51").Value = "'0.05943"
$ws.Range("E51").Value = "  +0.39%  "
